{"js": "const body = context.document.body;\n\n// 1. Merge the split \"Versi\" + \"on\" runs into a single \"Version\" run.\n//    Office.js's search() finds \"Version\" even though it spans two runs;\n//    replacing it with itself coalesces the text into one run while\n//    preserving the surrounding proofErr spell-check markers.\nconst versionResults = body.search(\"Version\", { matchCase: true });\nversionResults.load(\"items\");\nawait context.sync();\nversionResults.items[0].insertText(\"Version\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Bump the version number: \" 2\" -> \" 1.\"\nconst numberResults = body.search(\" 2\", { matchCase: true });\nnumberResults.load(\"items\");\nawait context.sync();\nnumberResults.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Remove the now-redundant trailing \".\" run (the original period that\n//    used to follow the bookmark). After step 2 there are two \".\" matches\n//    in the paragraph (the one just inserted, and the original trailing\n//    one) -- we want the last (original) one removed.\nconst dotResults = body.search(\".\", { matchCase: true });\ndotResults.load(\"items\");\nawait context.sync();\ndotResults.items[dotResults.items.length - 1].delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Merge the split \"Versi\" + \"on\" runs into a single \"Version\" run.\n#    Find.Execute matches \"Version\" even though it spans two runs; replacing\n#    it with itself coalesces the text into one run while the surrounding\n#    proofErr spell-check markers stay where they are.\n$r1 = $d.Content\n$r1.Find.Execute(\"Version\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2)\n\n# 2. Remove the trailing \".\" run (the period that used to follow the\n#    bookmark). At this point in the document there is exactly one \".\"\n#    so this is unambiguous; doing it now (before changing \"2\" to \"1.\")\n#    avoids any confusion with the period we are about to introduce below.\n$r2 = $d.Content\n$r2.Find.Execute(\".\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$r2.Delete()\n\n# 3. Bump the version number: \" 2\" -> \" 1.\"\n$r3 = $d.Content\n$r3.Find.Execute(\" 2\", $false, $false, $false, $false, $false, $true, 1, $false, \" 1.\", 2)\n"}
